# Weekly update: new "Femacal de La Calera - Sandia" prices were added.
# Three new observation rows (Extra / Primera / Segunda, Region de O'Higgins,
# fecha 2023-02-08 / serial 44984) are inserted above the existing history,
# which shifts the rest of the table down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 608; existing rows 608:699 become 611:702.
$ws.Rows("608:610").Insert()

# Row 608 - Extra
$ws.Cells.Item(608, 1).Value = 3
$ws.Cells.Item(608, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(608, 3).Value = "Coquimbo"
$ws.Cells.Item(608, 4).Value = 44984
$ws.Cells.Item(608, 5).Value = 5
$ws.Cells.Item(608, 6).Value = 100112028
$ws.Cells.Item(608, 7).Value = "Sandia"
$ws.Cells.Item(608, 8).Value = "Sin especificar"
$ws.Cells.Item(608, 9).Value = "Extra"
$ws.Cells.Item(608, 10).Value = 780
$ws.Cells.Item(608, 11).Value = 2800
$ws.Cells.Item(608, 12).Value = 3000
$ws.Cells.Item(608, 13).Value = 2897
$ws.Cells.Item(608, 14).Value = "$/unidad"
$ws.Cells.Item(608, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(608, 16).Value = 2897
$ws.Cells.Item(608, 17).Value = 1
$ws.Cells.Item(608, 18).Value = "Hortaliza"

# Row 609 - Primera
$ws.Cells.Item(609, 1).Value = 3
$ws.Cells.Item(609, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(609, 3).Value = "Coquimbo"
$ws.Cells.Item(609, 4).Value = 44984
$ws.Cells.Item(609, 5).Value = 5
$ws.Cells.Item(609, 6).Value = 100112028
$ws.Cells.Item(609, 7).Value = "Sandia"
$ws.Cells.Item(609, 8).Value = "Sin especificar"
$ws.Cells.Item(609, 9).Value = "Primera"
$ws.Cells.Item(609, 10).Value = 4150
$ws.Cells.Item(609, 11).Value = 2000
$ws.Cells.Item(609, 12).Value = 2200
$ws.Cells.Item(609, 13).Value = 2183
$ws.Cells.Item(609, 14).Value = "$/unidad"
$ws.Cells.Item(609, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(609, 16).Value = 2183
$ws.Cells.Item(609, 17).Value = 1
$ws.Cells.Item(609, 18).Value = "Hortaliza"

# Row 610 - Segunda
$ws.Cells.Item(610, 1).Value = 3
$ws.Cells.Item(610, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(610, 3).Value = "Coquimbo"
$ws.Cells.Item(610, 4).Value = 44984
$ws.Cells.Item(610, 5).Value = 5
$ws.Cells.Item(610, 6).Value = 100112028
$ws.Cells.Item(610, 7).Value = "Sandia"
$ws.Cells.Item(610, 8).Value = "Sin especificar"
$ws.Cells.Item(610, 9).Value = "Segunda"
$ws.Cells.Item(610, 10).Value = 2660
$ws.Cells.Item(610, 11).Value = 1500
$ws.Cells.Item(610, 12).Value = 1700
$ws.Cells.Item(610, 13).Value = 1673
$ws.Cells.Item(610, 14).Value = "$/unidad"
$ws.Cells.Item(610, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(610, 16).Value = 1673
$ws.Cells.Item(610, 17).Value = 1
$ws.Cells.Item(610, 18).Value = "Hortaliza"
